$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report record was added for "Poroto granado" at Vega Modelo de
# Temuco. It lands right above the existing row that used to be row 36, so
# insert a fresh row there; this pushes the former rows 36:90 down to 37:91
# (Excel carries the surrounding formatting/number-format along for us).
$ws.Rows("36:36").Insert()

# Fill in the newly inserted row 36 with the new record's values.
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 44935
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 100112030
$ws.Range("G36").Value = "Poroto granado"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 55
$ws.Range("K36").Value = 50000
$ws.Range("L36").Value = 50000
$ws.Range("M36").Value = 50000
$ws.Range("N36").Value = "$/saco 25 kilos"
$ws.Range("O36").Value = "Región del Maule"
$ws.Range("P36").Value = 2000
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
